# Append the new daily-log row (row 35) to each of the 4 data sheets.
# Each sheet already has data rows 2..34; the new row duplicates the
# previous (last) row's values except for the timestamp column (A),
# which advances to the next day, and (for one sheet) the D/H columns
# which carry the freshly-measured "actual length" values.

$wb = $excel.ActiveWorkbook

$newRow = 35
$prevRow = 34

# Per-sheet column values for the new row. Columns B, C, D, E are text
# (hex byte strings); A, F, G, H, I are numeric.
$rowsData = @{
    "DE_LFT_#1" = @{
        A = "45821.43564814814"
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
        D = "0x01,0x74"
        E = "0x14"
        F = 380
        G = "7.598631275147109e+23"
        H = 372
        I = 14
    }
    "DE_LFT_#2" = @{
        A = "45821.43564814814"
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
        D = "0x01,0x74"
        E = "0xe"
        F = 380
        G = "5.68432987514711e+23"
        H = 372
        I = 14
    }
    "DE_PLT_#1" = @{
        A = "45821.43564814814"
        B = "0x00,0x82"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x00,0x80"
        E = "0x7"
        F = 130
        G = "5.68631262647114e+23"
        H = 128
        I = 7
    }
    "DE_PLT_#2" = @{
        A = "45821.43564814814"
        B = "0x00,0x82"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x00,0x80"
        E = "0x3"
        F = 130
        G = "9.85046333984776e+23"
        H = 128
        I = 3
    }
}

foreach ($ws in $wb.Worksheets) {
    $data = $rowsData[$ws.Name]
    if ($data -eq $null) { continue }

    # Timestamp column keeps the same number format/style as the row above it.
    $ws.Cells.Item($newRow, 1).Value = [double]$data.A
    $ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($prevRow, 1).NumberFormat

    $ws.Cells.Item($newRow, 2).Value = $data.B
    $ws.Cells.Item($newRow, 3).Value = $data.C
    $ws.Cells.Item($newRow, 4).Value = $data.D
    $ws.Cells.Item($newRow, 5).Value = $data.E

    $ws.Cells.Item($newRow, 6).Value = $data.F
    $ws.Cells.Item($newRow, 7).Value = [double]$data.G
    $ws.Cells.Item($newRow, 8).Value = $data.H
    $ws.Cells.Item($newRow, 9).Value = $data.I
}
